$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.268.14'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.920.69'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -1.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.51'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("E6").Value = '  -1.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4854'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3856'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07421'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9522'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.98'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07810'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.914.05'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.556'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.667'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.29'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008905'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.263.73'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.05'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.177'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("E23").Value = '  +2.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.16%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.933'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.53%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.64'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.127'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.07'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.039'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08914'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.365'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("E33").Value = '  +4.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7822'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.689'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.794'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.43%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02058'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.40%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.131'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05384'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.63%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5592'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.035'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.132'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.602'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1539'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.83'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.02%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4936'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.69'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.28%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.48%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.680'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.52'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +4.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06153'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.05%  '
